$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H19").Value = 259.8
$ws.Range("I19").Value = 251
$ws.Range("J19").Value = 273
$ws.Range("K19").Value = 251
$ws.Range("L19").Value = 273
$ws.Range("M19").Value = -76
$ws.Range("N19").Value = -623
$ws.Range("H33").Value = 357.65216
$ws.Range("I33").Value = 356.73685
$ws.Range("J33").Value = 362
$ws.Range("K33").Value = 356.73685
$ws.Range("L33").Value = 362
$ws.Range("M33").Value = -127.73685
$ws.Range("N33").Value = -820
$ws.Range("H137").Value = 970797.25
$ws.Range("I137").Value = 556412.5
$ws.Range("J137").Value = 1385182
$ws.Range("K137").Value = 1669237.5
$ws.Range("L137").Value = 4155546
$ws.Range("M137").Value = -1666687.5
$ws.Range("N137").Value = -4160646
$ws.Range("H138").Value = 3854.4167
$ws.Range("J138").Value = 5379.6665
$ws.Range("L138").Value = 16138.9995
$ws.Range("N138").Value = -26418.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1157331
$ws.Range("I61").Value = 1242837
$ws.Range("K61").Value = 1242837
$ws.Range("M61").Value = -1242625
$ws.Range("H74").Value = 2909770.8
$ws.Range("I74").Value = 3474216.5
$ws.Range("K74").Value = 3474216.5
$ws.Range("M74").Value = -3473342.5
$ws.Range("H77").Value = 2909770.8
$ws.Range("I77").Value = 3474216.5
$ws.Range("K77").Value = 17371082.5
$ws.Range("M77").Value = -17366714.5
$ws.Range("H97").Value = 686.6667
$ws.Range("I97").Value = 694.63635
$ws.Range("K97").Value = 694.63635
$ws.Range("M97").Value = -198.63635
$ws.Range("H122").Value = 2578.5833
$ws.Range("I122").Value = 1524.375
$ws.Range("J122").Value = 4687
$ws.Range("K122").Value = 4573.125
$ws.Range("L122").Value = 14061
$ws.Range("M122").Value = -2123.125
$ws.Range("N122").Value = -18961
$ws.Range("H132").Value = 730317.5600000001
$ws.Range("I132").Value = 1013132.25
$ws.Range("J132").Value = 7568.8887
$ws.Range("K132").Value = 3039396.75
$ws.Range("L132").Value = 22706.6661
$ws.Range("M132").Value = -3036866.75
$ws.Range("N132").Value = -27766.6661
$ws.Range("H136").Value = 1157331
$ws.Range("I136").Value = 1242837
$ws.Range("K136").Value = 3728511
$ws.Range("M136").Value = -3725961

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 520443.47
$ws.Range("I134").Value = 637753.7
$ws.Range("K134").Value = 1913261.1
$ws.Range("M134").Value = -1910726.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 20000
$ws.Range("J18").Value = 20000
$ws.Range("L18").Value = 20000
$ws.Range("N18").Value = -20460
$ws.Range("H22").Value = 844.2083
$ws.Range("I22").Value = 848.9375
$ws.Range("J22").Value = 834.75
$ws.Range("K22").Value = 848.9375
$ws.Range("L22").Value = 834.75
$ws.Range("M22").Value = -498.9375
$ws.Range("N22").Value = -1534.75
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H96").Value = 39871.332
$ws.Range("J96").Value = 39871.332
$ws.Range("L96").Value = 39871.332
$ws.Range("N96").Value = -45363.332
$ws.Range("H107").Value = 866.1
$ws.Range("I107").Value = 1174.8334
$ws.Range("J107").Value = 403
$ws.Range("K107").Value = 1174.8334
$ws.Range("L107").Value = 403
$ws.Range("M107").Value = 745.1666
$ws.Range("N107").Value = -4243

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1174.4
$ws.Range("J92").Value = 1733.1666
$ws.Range("L92").Value = 5199.4998
$ws.Range("N92").Value = -7695.4998
$ws.Range("H114").Value = 1722.5238
$ws.Range("I114").Value = 94.21429000000001
$ws.Range("J114").Value = 4979.143
$ws.Range("K114").Value = 282.64287
$ws.Range("L114").Value = 14937.429
$ws.Range("M114").Value = 2971.35713
$ws.Range("N114").Value = -21445.429
$ws.Range("H121").Value = 10527462
$ws.Range("J121").Value = 1450.8572
$ws.Range("L121").Value = 4352.571599999999
$ws.Range("N121").Value = -6972.571599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 23380
$ws.Range("J39").Value = 23380
$ws.Range("L39").Value = 23380
$ws.Range("N39").Value = -24444
$ws.Range("H102").Value = 2522.4333
$ws.Range("I102").Value = 1498.4546
$ws.Range("J102").Value = 5338.375
$ws.Range("K102").Value = 1498.4546
$ws.Range("L102").Value = 5338.375
$ws.Range("M102").Value = 123.5454
$ws.Range("N102").Value = -8582.375
$ws.Range("H132").Value = 247410.19
$ws.Range("J132").Value = 3689.6667
$ws.Range("L132").Value = 11069.0001
$ws.Range("N132").Value = -16129.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 895.6875
$ws.Range("J55").Value = 1610
$ws.Range("L55").Value = 1610
$ws.Range("N55").Value = -1956
$ws.Range("H59").Value = 42500.5
$ws.Range("J59").Value = 42500.5
$ws.Range("L59").Value = 42500.5
$ws.Range("N59").Value = -43808.5
$ws.Range("H61").Value = 3503.1875
$ws.Range("I61").Value = 1983.5
$ws.Range("K61").Value = 1983.5
$ws.Range("M61").Value = -1781.5
$ws.Range("H70").Value = 11500
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 11500
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 11500
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -12040
$ws.Range("H73").Value = 11500
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 11500
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 11500
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -13372
$ws.Range("H82").Value = 1737.6666
$ws.Range("I82").Value = 1370.3077
$ws.Range("K82").Value = 1370.3077
$ws.Range("M82").Value = -1009.3077
$ws.Range("H85").Value = 1737.6666
$ws.Range("I85").Value = 1370.3077
$ws.Range("K85").Value = 1370.3077
$ws.Range("M85").Value = -122.3077000000001
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H93").Value = 1789.55
$ws.Range("I93").Value = 1335.8182
$ws.Range("K93").Value = 1335.8182
$ws.Range("M93").Value = -87.81819999999993
$ws.Range("H111").Value = 30351
$ws.Range("I111").Value = 30351
$ws.Range("K111").Value = 30351
$ws.Range("M111").Value = -26261
$ws.Range("H113").Value = 3503.1875
$ws.Range("I113").Value = 1983.5
$ws.Range("K113").Value = 1983.5
$ws.Range("M113").Value = 186.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 30000
$ws.Range("J56").Value = 30000
$ws.Range("L56").Value = 30000
$ws.Range("N56").Value = -31428
$ws.Range("H114").Value = 80397.8
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 80397.8
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 80397.8
$ws.Range("M114").ClearContents()
$ws.Range("N114").Value = -89075.8
$ws.Range("H136").Value = 9297706
$ws.Range("I136").Value = 9774127
$ws.Range("K136").Value = 29322381
$ws.Range("M136").Value = -29319831
